# feat: add 2022-Q1 data
#
# Before:  2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计
# After:   2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计
#
# The old "总计" sheet (sheetId 6) is repurposed to hold the new "2022-Q1"
# fund-detail breakdown (matching the format already used by "2021-Q4"),
# and a brand-new "总计" sheet is appended after it, carrying the old
# totals table plus a freshly prepended 2022-Q1 summary row.
#
# NOTE: sheet-object variables in this host are positional handles, not
# stable references -- Worksheets.Add()/Move() renumber everything and
# silently invalidate variables captured beforehand. So: finish ALL
# structural changes (rename/add/move) first, THEN re-fetch every sheet
# by name and only afterwards touch cell contents.

$wb = $excel.ActiveWorkbook

# --- structural changes only, in this block ---------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$wb.Worksheets.Add() | Out-Null
$newSheet = $wb.Worksheets.Item("Sheet1")
$q1ForAnchor = $wb.Worksheets.Item("2022-Q1")
$newSheet.Move($null, $q1ForAnchor)

$newSheet2 = $wb.Worksheets.Item("Sheet1")
$newSheet2.Name = "总计"
# ------------------------------------------------------------------------

# From here on, no more Add()/Move()/Delete() calls -- safe to keep refs.
$styleDonor = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Item("2022-Q1")
$total = $wb.Worksheets.Item("总计")

# "2022-Q1" reused the old "总计" sheet, which had 6 rows of unrelated
# data -- wipe it clean before writing the new fund-detail table.
$q1.Cells.Clear()

# ---------------------------------------------------------------------
# "2022-Q1": fund-level detail (same shape as the existing "2021-Q4").
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $donor = $styleDonor.Cells.Item(1, $col)
    $dest = $q1.Cells.Item(1, $col)
    $donor.Copy($dest)
    $dest.Value = $headers[$col - 2]
}

# Data rows: index(A), code(B), name(C), size(D), position(E), ratio(F), value(G), rank(H)
# B/D/E/F/G hold digit-only strings in the source file (e.g. fund code
# "161838", ratio "5.92") -- force Text format so they land as literal
# strings instead of being auto-coerced to numbers.
$rows = @(
    @(0, "161838", "银华创业板两年定期开放混合", "10.44", "95.40", "5.92", "0.6180", 5),
    @(1, "159851", "华宝中证金融科技主题ETF", "3.16", "98.58", "2.92", "0.0923", 9),
    @(2, "516100", "华夏中证金融科技主题交易型开放式指数证券投资基金", "0.68", "96.91", "2.90", "0.0197", 9)
)

$r = 2
foreach ($row in $rows) {
    $styleDonor.Cells.Item($r, 1).Copy($q1.Cells.Item($r, 1))
    $q1.Cells.Item($r, 1).Value = $row[0]

    $q1.Cells.Item($r, 2).NumberFormat = "@"
    $q1.Cells.Item($r, 2).Value = $row[1]

    $q1.Cells.Item($r, 3).Value = $row[2]

    $q1.Cells.Item($r, 4).NumberFormat = "@"
    $q1.Cells.Item($r, 4).Value = $row[3]

    $q1.Cells.Item($r, 5).NumberFormat = "@"
    $q1.Cells.Item($r, 5).Value = $row[4]

    $q1.Cells.Item($r, 6).NumberFormat = "@"
    $q1.Cells.Item($r, 6).Value = $row[5]

    $q1.Cells.Item($r, 7).NumberFormat = "@"
    $q1.Cells.Item($r, 7).Value = $row[6]

    $q1.Cells.Item($r, 8).Value = $row[7]

    $r++
}

# ---------------------------------------------------------------------
# "总计": old totals table, plus a prepended 2022-Q1 summary row.
# ---------------------------------------------------------------------
$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($col = 2; $col -le 4; $col++) {
    $donor = $styleDonor.Cells.Item(1, $col)
    $dest = $total.Cells.Item(1, $col)
    $donor.Copy($dest)
    $dest.Value = $totalHeaders[$col - 2]
}

# Totals rows: index(A), date(B), count(C), value(D) -- newest quarter on top.
$totalRows = @(
    @(0, "2022-Q1", 3, 0.73),
    @(1, "2021-Q4", 6, 0.84),
    @(2, "2021-Q3", 4, 1.63),
    @(3, "2021-Q2", 6, 4.09),
    @(4, "2021-Q1", 6, 3.39),
    @(5, "2020-Q4", 4, 3.59)
)

$r = 2
foreach ($row in $totalRows) {
    $styleDonor.Cells.Item($r, 1).Copy($total.Cells.Item($r, 1))
    $total.Cells.Item($r, 1).Value = $row[0]

    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]

    $r++
}
